# ManageInvoiceTestData.xlsx -- "Test and bat update"
#
# The test-data row on Sheet1 (row 2) is refreshed with the next batch of
# automation values: OrderDate (column A) and OverageID (column L) move on
# to the next recorded order/overage pair. Everything else in the row is
# left alone.
#
# Values are written through a Formula/"paste values" round-trip rather
# than a plain `.Value =` assignment so that date-shaped text such as
# "11-18-2021" is stored as literal text (matching the workbook's existing
# General-formatted, shared-string cells) instead of being auto-parsed
# into a date serial number, and so the cell keeps its existing
# border/fill formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A2: OrderDate -> "11-18-2021"
$ws.Range("A2").Formula = "=""11-18-2021"""
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

# L2: OverageID -> "58233443"
$ws.Range("L2").Formula = "=""58233443"""
$ws.Range("L2").Copy()
$ws.Range("L2").PasteSpecial(-4163)

$excel.CutCopyMode = 0
